$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 833.3333
$ws.Range("I12").Value = 799
$ws.Range("J12").Value = 850.5
$ws.Range("K12").Value = 799
$ws.Range("L12").Value = 850.5
$ws.Range("M12").Value = -629
$ws.Range("N12").Value = -1190.5
$ws.Range("H19").Value = 5609.2383
$ws.Range("J19").Value = 5524.8335
$ws.Range("L19").Value = 5524.8335
$ws.Range("N19").Value = -5874.8335
$ws.Range("H32").Value = 1450.3334
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1450.3334
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1450.3334
$ws.Range("M32").ClearContents() | Out-Null
$ws.Range("N32").Value = -2102.3334
$ws.Range("H40").Value = 33333332
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 33333332
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 33333332
$ws.Range("M40").ClearContents() | Out-Null
$ws.Range("N40").Value = -33333682
$ws.Range("H43").Value = 262114.56
$ws.Range("I43").Value = 6583.7
$ws.Range("J43").Value = 687999.3
$ws.Range("K43").Value = 6583.7
$ws.Range("L43").Value = 687999.3
$ws.Range("M43").Value = -6514.7
$ws.Range("N43").Value = -688137.3
$ws.Range("H64").Value = 66674570
$ws.Range("I64").Value = 125006730
$ws.Range("K64").Value = 125006730
$ws.Range("M64").Value = -125006482
$ws.Range("H67").Value = 66674570
$ws.Range("I67").Value = 125006730
$ws.Range("K67").Value = 125006730
$ws.Range("M67").Value = -125005872
$ws.Range("H98").Value = 2929.2856
$ws.Range("I98").Value = 3127.5
$ws.Range("K98").Value = 3127.5
$ws.Range("M98").Value = -1629.5
$ws.Range("H116").Value = 19239822
$ws.Range("I116").Value = 125004250
$ws.Range("K116").Value = 125004250
$ws.Range("M116").Value = -125000808
$ws.Range("H122").Value = 2929.2856
$ws.Range("I122").Value = 3127.5
$ws.Range("K122").Value = 9382.5
$ws.Range("M122").Value = -6932.5
$ws.Range("H137").Value = 3452.8
$ws.Range("I137").Value = 3300.6667
$ws.Range("K137").Value = 9902.000100000001
$ws.Range("M137").Value = -7352.000100000001
$ws.Range("H138").Value = 1617157
$ws.Range("I138").Value = 1337.0385
$ws.Range("K138").Value = 4011.1155
$ws.Range("M138").Value = 1128.8845

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 125005470
$ws.Range("I2").Value = 1606.3334
$ws.Range("K2").Value = 1606.3334
$ws.Range("M2").Value = -1493.3334
$ws.Range("H32").Value = 1671327.5
$ws.Range("I32").Value = 1739577
$ws.Range("K32").Value = 1739577
$ws.Range("M32").Value = -1739290
$ws.Range("H74").Value = 95863.17999999999
$ws.Range("I74").Value = 203899
$ws.Range("K74").Value = 203899
$ws.Range("M74").Value = -203025
$ws.Range("H77").Value = 95863.17999999999
$ws.Range("I77").Value = 203899
$ws.Range("K77").Value = 1019495
$ws.Range("M77").Value = -1015127
$ws.Range("H102").Value = 1572.5
$ws.Range("I102").Value = 1572.5
$ws.Range("K102").Value = 1572.5
$ws.Range("M102").Value = 49.5
$ws.Range("H116").Value = 125005470
$ws.Range("I116").Value = 1606.3334
$ws.Range("K116").Value = 1606.3334
$ws.Range("M116").Value = 687.6666
$ws.Range("H132").Value = 3273.8408
$ws.Range("I132").Value = 1458.2424
$ws.Range("J132").Value = 8720.637000000001
$ws.Range("K132").Value = 4374.7272
$ws.Range("L132").Value = 26161.911
$ws.Range("M132").Value = -1844.7272
$ws.Range("N132").Value = -31221.911

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 125005470
$ws.Range("I3").Value = 1606.3334
$ws.Range("K3").Value = 1606.3334
$ws.Range("M3").Value = -1492.3334
$ws.Range("H22").Value = 286.18182
$ws.Range("I22").Value = 280.83334
$ws.Range("J22").Value = 292.6
$ws.Range("K22").Value = 280.83334
$ws.Range("L22").Value = 292.6
$ws.Range("M22").Value = -107.83334
$ws.Range("N22").Value = -638.6
$ws.Range("H86").Value = 6276397
$ws.Range("J86").Value = 2023.3572
$ws.Range("L86").Value = 2023.3572
$ws.Range("N86").Value = -4269.3572
$ws.Range("H89").Value = 6276397
$ws.Range("J89").Value = 2023.3572
$ws.Range("L89").Value = 10116.786
$ws.Range("N89").Value = -21348.786
$ws.Range("H99").Value = 4331408.5
$ws.Range("J99").Value = 7577485
$ws.Range("L99").Value = 7577485
$ws.Range("N99").Value = -7580481
$ws.Range("H134").Value = 4552.6826
$ws.Range("I134").Value = 1677.186
$ws.Range("K134").Value = 5031.558
$ws.Range("M134").Value = -2496.558

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 55.7
$ws.Range("I7").Value = 56.15
$ws.Range("K7").Value = 56.15
$ws.Range("M7").Value = 56.85
$ws.Range("H16").Value = 7267.7334
$ws.Range("I16").Value = 6229
$ws.Range("K16").Value = 6229
$ws.Range("M16").Value = -5942
$ws.Range("H22").Value = 846.2857
$ws.Range("J22").Value = 699
$ws.Range("L22").Value = 699
$ws.Range("N22").Value = -1399
$ws.Range("H58").Value = 5792.731
$ws.Range("I58").Value = 2382.2122
$ws.Range("J58").Value = 11716.263
$ws.Range("K58").Value = 2382.2122
$ws.Range("L58").Value = 11716.263
$ws.Range("M58").Value = -2179.2122
$ws.Range("N58").Value = -12122.263
$ws.Range("H105").Value = 3574145.2
$ws.Range("I105").Value = 3969939.2
$ws.Range("J105").Value = 12000
$ws.Range("K105").Value = 3969939.2
$ws.Range("L105").Value = 12000
$ws.Range("M105").Value = -3968192.2
$ws.Range("N105").Value = -15494
$ws.Range("H113").Value = 7267.7334
$ws.Range("I113").Value = 6229
$ws.Range("K113").Value = 6229
$ws.Range("M113").Value = -4059
$ws.Range("H132").Value = 6267.095
$ws.Range("I132").Value = 3595.3333
$ws.Range("J132").Value = 9829.444
$ws.Range("K132").Value = 10785.9999
$ws.Range("L132").Value = 29488.332
$ws.Range("M132").Value = -8255.999899999999
$ws.Range("N132").Value = -34548.33199999999
$ws.Range("H134").Value = 7729.4517
$ws.Range("I134").Value = 3653.25
$ws.Range("K134").Value = 10959.75
$ws.Range("M134").Value = -8424.75
$ws.Range("H136").Value = 5792.731
$ws.Range("I136").Value = 2382.2122
$ws.Range("J136").Value = 11716.263
$ws.Range("K136").Value = 7146.6366
$ws.Range("L136").Value = 35148.789
$ws.Range("M136").Value = -4596.6366
$ws.Range("N136").Value = -40248.789

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4600
$ws.Range("J104").Value = 5000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -20242
$ws.Range("H122").Value = 1179730.2
$ws.Range("I122").Value = 2358238.2
$ws.Range("J122").Value = 1222.1666
$ws.Range("K122").Value = 21224143.8
$ws.Range("L122").Value = 10999.4994
$ws.Range("M122").Value = -21221693.8
$ws.Range("N122").Value = -15899.4994
$ws.Range("H129").Value = 9315451
$ws.Range("I129").Value = 422.3846
$ws.Range("J129").Value = 33534524
$ws.Range("K129").Value = 1267.1538
$ws.Range("L129").Value = 100603572
$ws.Range("M129").Value = 3732.8462
$ws.Range("N129").Value = -100613572
$ws.Range("H131").Value = 2172.0188
$ws.Range("I131").Value = 1310.125
$ws.Range("J131").Value = 2325.2444
$ws.Range("K131").Value = 3930.375
$ws.Range("L131").Value = 6975.733200000001
$ws.Range("M131").Value = 1109.625
$ws.Range("N131").Value = -17055.7332

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17475.477
$ws.Range("I132").Value = 16686.625
$ws.Range("J132").Value = 19999.8
$ws.Range("K132").Value = 50059.875
$ws.Range("L132").Value = 59999.39999999999
$ws.Range("M132").Value = -47529.875
$ws.Range("N132").Value = -65059.39999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 544.75
$ws.Range("J16").Value = 450
$ws.Range("L16").Value = 450
$ws.Range("N16").Value = -790
$ws.Range("H46").Value = 1502809.8
$ws.Range("I46").Value = 4311660
$ws.Range("J46").Value = 4756.3335
$ws.Range("K46").Value = 4311660
$ws.Range("L46").Value = 4756.3335
$ws.Range("M46").Value = -4311472
$ws.Range("N46").Value = -5132.3335
$ws.Range("H55").Value = 533
$ws.Range("I55").Value = 250
$ws.Range("J55").Value = 568.375
$ws.Range("K55").Value = 250
$ws.Range("L55").Value = 568.375
$ws.Range("M55").Value = -77
$ws.Range("N55").Value = -914.375
$ws.Range("H100").Value = 4106.467
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 4106.467
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 4106.467
$ws.Range("M100").ClearContents() | Out-Null
$ws.Range("N100").Value = -5188.467
$ws.Range("H136").Value = 11915
$ws.Range("I136").Value = 4821.25
$ws.Range("K136").Value = 14463.75
$ws.Range("M136").Value = -11913.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19609148
$ws.Range("I107").Value = 1320.7142
$ws.Range("K107").Value = 3962.1426
$ws.Range("M107").Value = -2042.1426
$ws.Range("H126").Value = 4185.5713
$ws.Range("I126").Value = 3659.8
$ws.Range("K126").Value = 10979.4
$ws.Range("M126").Value = -8509.400000000001
$ws.Range("H132").Value = 12197535
$ws.Range("I132").Value = 13515687
$ws.Range("J132").Value = 4636.25
$ws.Range("K132").Value = 40547061
$ws.Range("L132").Value = 13908.75
$ws.Range("M132").Value = -40544531
$ws.Range("N132").Value = -18968.75
$ws.Range("H136").Value = 20023598
$ws.Range("I136").Value = 33334436
$ws.Range("J136").Value = 57339.9
$ws.Range("K136").Value = 100003308
$ws.Range("L136").Value = 172019.7
$ws.Range("M136").Value = -100000758
$ws.Range("N136").Value = -177119.7
